# Insert a new data row at row 205 (pushing existing rows 205-219 down to
# 206-220), then populate the newly inserted row with its values.
# This mirrors the XML diff: a new record (Berenjena, Primera) was added
# to the "Hortaliza ... Berenjena" sheet, causing the dimension to grow
# from A1:R219 to A1:R220.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 205..219 down to 206..220, creating a blank row 205.
$ws.Rows(205).Insert()

# Populate the new row 205 with the inserted record's data.
$ws.Range("A205").Value = 6
$ws.Range("B205").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C205").Value = "Metropolitana"
$ws.Range("D205").Value = 44746
$ws.Range("E205").Value = 13
$ws.Range("F205").Value = 100112001
$ws.Range("G205").Value = "Berenjena"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 180
$ws.Range("K205").Value = 10000
$ws.Range("L205").Value = 12000
$ws.Range("M205").Value = 10889
$ws.Range("N205").Value = "$/caja 50 unidades"
$ws.Range("O205").Value = "Región de Arica y Parinacota"
$ws.Range("P205").Value = 218
$ws.Range("Q205").Value = 50
$ws.Range("R205").Value = "Hortaliza"

# Ensure the date cell keeps the same date number format as the rest of
# column D (style index 2 in the original workbook).
$ws.Range("D205").NumberFormat = $ws.Range("D206").NumberFormat
